# TestSuite.xlsx template update
#
# The "Login and Registe" test-suite row is renamed to "Login module" and
# the sheet is trimmed down to just the header row plus that single data
# row (the other sample rows - Record Modify Test / TS_Add new record /
# TS_Delete a record / the bare "No" / "Yes" rows - are removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the module/test-case name in row 2.
$ws.Range("B2").Value = "Login module"

# Drop rows 3-6; the used range collapses back down to A1:F2.
$ws.Rows("3:6").Delete()

# Reset the active selection away from the old F2 hyperlink cell.
$ws.Range("A1").Select()
